$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.638.65'
$ws.Range('E2').Value = '  -2.88%  '
$ws.Range('D3').Value = '3.479.39'
$ws.Range('E3').Value = '  -1.88%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '569.78'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.03%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '176.10'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -8.87%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.630'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +3.11%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.625'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.85%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.154'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.86%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '53.68'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -5.48%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000269'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.70%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '9.10'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -4.75%  '
$ws.Range('D14').Value = '4.040.68'
$ws.Range('E14').Value = '  -1.79%  '
$ws.Range('D15').Value = '3.489.17'
$ws.Range('E15').Value = '  -1.68%  '
$ws.Range('E16').Value = '  -0.51%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '18.12'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.82%  '
$ws.Range('D18').Value = '65.674.86'
$ws.Range('E18').Value = '  -3.00%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.98'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.03%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.992'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.46%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '410.72'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.43%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.16'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +4.11%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '84.71'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.38%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '4.19'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.57%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '12.58'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +3.99%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '10.78'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -5.33%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.80'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -5.99%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.89'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.38%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '29.99'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.17%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '617.42'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -10.41%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.30'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -8.16%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '11.53'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.83%  '
$ws.Range('E33').Value = '  -3.89%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '59.10'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.81%  '
$ws.Range('E35').Value = '  +7.42%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.00'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.06%  '
$ws.Range('D37').Value = '0.0₃0786'
$ws.Range('E37').Value = '  -5.16%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '36.77'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -6.36%  '
$ws.Range('D39').Value = '3.297.93'
$ws.Range('E39').Value = '  +8.53%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.375'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -6.03%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.33'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.93%  '
$ws.Range('E42').Value = '  +0.11%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.86'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -5.79%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.24'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.19%  '
$ws.Range('B45').Value = 'VeChain'
$ws.Range('C45').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0412'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.60%  '
$ws.Range('B46').Value = 'Fetch.AI'
$ws.Range('C46').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.48'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -7.84%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.69'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.87%  '
$ws.Range('E48').Value = '  -0.35%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '138.88'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.27%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.34'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -8.60%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.79'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +5.29%  '
